# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps produced by a fresh handback report run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (column G) for the first data row
$overview.Range("G2").Value = "2016-09-01 05:09:29"

# zh-cn sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K)
$zhcn.Range("H2").Value = "2016-09-01 05:09:24"
$zhcn.Range("K2").Value = "2016-09-01 05:09:41"

# de-de sheet: Correspond Handback DateTime (K) -- H2 shares the value with Overview!G2
$dede.Range("K2").Value = "2016-09-01 05:09:48"
